$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume/Coin/Link columns keep their original plain-text
# formatting instead of Excel auto-coercing numeric-looking strings (e.g.
# "219.30", "1.00", "0.0000256") into actual numbers.
$ws.Range("B2:E51").NumberFormat = "@"

# --- Column D (Price) updates for existing rows ---
$dUpdates = @{
    2  = "91.229.47"
    3  = "3.102.08"
    5  = "219.30"
    6  = "618.58"
    10 = "3.099.34"
    11 = "0.696"
    13 = "0.0000256"
    15 = "91.074.17"
    16 = "33.12"
    17 = "3.673.03"
    18 = "3.082.07"
    20 = "0.0000233"
    21 = "13.84"
    22 = "431.29"
    23 = "8.55"
    25 = "5.58"
    26 = "11.92"
    27 = "83.74"
    28 = "3.263.06"
    31 = "1.00"
    32 = "8.71"
    33 = "3.91"
    34 = "520.03"
    35 = "7.01"
    36 = "0.142"
    39 = "23.02"
    40 = "22.31"
    47 = "142.35"
}

foreach ($row in $dUpdates.Keys) {
    $ws.Range("D$row").Value = $dUpdates[$row]
}

# --- Column E (Volume 1h) updates for existing rows ---
$eUpdates = @{
    2  = "  +4.71%  "
    3  = "  +1.54%  "
    5  = "  +5.25%  "
    6  = "  +0.21%  "
    7  = "  +5.80%  "
    8  = "  +16.23%  "
    9  = "  -0.03%  "
    10 = "  +1.67%  "
    11 = "  +20.41%  "
    12 = "  +7.11%  "
    13 = "  +9.51%  "
    14 = "  +3.41%  "
    15 = "  +4.60%  "
    16 = "  +6.92%  "
    17 = "  +1.33%  "
    18 = "  +0.01%  "
    19 = "  +8.72%  "
    20 = "  +13.00%  "
    21 = "  +6.87%  "
    22 = "  +4.39%  "
    23 = "  +5.02%  "
    24 = "  +8.77%  "
    25 = "  +5.74%  "
    26 = "  +7.64%  "
    27 = "  +2.47%  "
    28 = "  +1.07%  "
    29 = "  -0.11%  "
    30 = "  +13.17%  "
    31 = "  -0.10%  "
    32 = "  +9.84%  "
    33 = "  +8.54%  "
    34 = "  +5.93%  "
    35 = "  +6.89%  "
    36 = "  +1.63%  "
    37 = "  +4.86%  "
    38 = "  +3.87%  "
    39 = "  +5.79%  "
    40 = "  +0.80%  "
    41 = "  -0.08%  "
    42 = "  +12.02%  "
    44 = "  +4.08%  "
    45 = "  +5.24%  "
    46 = "  +12.35%  "
    47 = "  -2.73%  "
    48 = "  +0.66%  "
}

foreach ($row in $eUpdates.Keys) {
    $ws.Range("E$row").Value = $eUpdates[$row]
}

# --- Rows 49-51: a new coin (ImmutableX) was inserted before FLOKI, shifting
#     FLOKI and Filecoin down by one row; the previous last row (Aave) drops off.
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").Value = "  +10.37%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "0.000260"
$ws.Range("E50").Value = "  +19.68%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "4.21"
$ws.Range("E51").Value = "  +9.42%  "
